$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.003982067108154
$ws.Range("B1").Value = 1.897334814071655
$ws.Range("C1").Value = 2.846660137176514
$ws.Range("D1").Value = 3.491264820098877
$ws.Range("E1").Value = 2.02833104133606
